# Fix properties for Tea_types
#
# The "properties" sheet (4th tab) had a stray "colOrder" row that needs to
# be removed, and that sheet should become the active/selected tab with the
# selection left on E7 (the cell below the last data row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("properties")

# Make "properties" the active sheet (this also updates the workbook's
# activeTab and moves tabSelected away from the previously active sheet).
$ws.Activate()

# Remove the "colOrder" row (row 2: Table | default | colOrder | string | [...])
$ws.Rows("2:2").Delete()

# Leave the selection on E7, matching the saved view state.
$ws.Range("E7").Select()
